# Update gh-pages output (合肥-漫展信息.xlsx): add the 2024.03.17 CW国潮
# "赵路内场" listing and refresh the attendance counters on the "展览"
# and "全部类型" sheets (they carry duplicate data).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Bump "想去人数" counters on the two untouched rows.
    $ws.Range("F2").Value = 1317
    $ws.Range("F3").Value = 1741

    # Insert a brand-new row 4 - this pushes the old rows 4-7 down to 5-8.
    $ws.Rows.Item(4).Insert()

    # The inserted row's left-most cell needs the same bordered/bold/centered
    # style used by the rest of column A; grab it from the row above.
    $ws.Range("A3").Copy()
    $ws.Range("A4").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # New row 4: 合肥·CW国潮动漫游戏嘉年华-赵路内场
    $ws.Range("A4").Value = 3
    # Format the date-shaped string as text first so Excel doesn't coerce
    # "2024.03.17" into a date serial, then drop back to the sheet's default
    # (unstyled) formatting to match the other rows' plain text cells.
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024.03.17"
    $ws.Range("B4").ClearFormats()
    $ws.Range("C4").Value = "合肥·CW国潮动漫游戏嘉年华-赵路内场"
    $ws.Range("D4").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Range("E4").Value = "2024.03.17 09:00-03.17 17:00"
    $ws.Range("F4").Value = 8
    $ws.Range("G4").Value = 238
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81954"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg"

    # Rows 5-8 (formerly 4-7) kept their own data via the shift; only the
    # running index in column A and the "想去人数" counters need refreshing.
    $ws.Range("A5").Value = 4
    $ws.Range("F5").Value = 69

    $ws.Range("A6").Value = 5
    $ws.Range("F6").Value = 6267

    $ws.Range("A7").Value = 6
    $ws.Range("F7").Value = 107

    $ws.Range("A8").Value = 7
}
